# Update the "Förändrad" (column C) date values for rows 2-18 from 2023-09-06
# (serial 45175) to 2023-09-14 (serial 45183), as produced by the automatic
# update of files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183   # Column C: "Förändrad"
}
